$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.311.21"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.932.55"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'251.19"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("D6").Value = "'0.7140"
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.3261"
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("D9").Value = "'27.34"
$ws.Range("E9").Value = "  +2.22%  "
$ws.Range("D10").Value = "'0.07193"
$ws.Range("E10").Value = "  +4.97%  "
$ws.Range("D11").Value = "'0.7998"
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("D12").Value = "'0.08091"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").Value = "1.928.36"
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").Value = "'5.434"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "'94.85"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").Value = "'14.83"
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("D17").Value = "30.324.28"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "'253.52"
$ws.Range("E18").Value = "  -4.17%  "
$ws.Range("D19").Value = "'0.000008094"
$ws.Range("E19").Value = "  +1.29%  "
$ws.Range("D20").Value = "'5.808"
$ws.Range("E20").Value = "  -0.99%  "
$ws.Range("D21").Value = "2.183.50"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "'6.917"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("D25").Value = "'9.710"
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("D26").Value = "'164.79"
$ws.Range("E26").Value = "  +2.78%  "
$ws.Range("D27").Value = "'19.28"
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("D28").Value = "'2.314"
$ws.Range("D29").Value = "'0.1279"
$ws.Range("E29").Value = "  -5.05%  "
$ws.Range("E30").Value = "  -0.40%  "
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").Value = "'4.207"
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("D34").Value = "'0.05208"
$ws.Range("E34").Value = "  +1.98%  "
$ws.Range("D35").Value = "'1.269"
$ws.Range("E35").Value = "  +4.81%  "
$ws.Range("D36").Value = "'0.7504"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").Value = "'2.770"
$ws.Range("D38").Value = "'0.01965"
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("D39").Value = "'2.803"
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("D40").Value = "'78.98"
$ws.Range("E40").Value = "  -2.95%  "
$ws.Range("D41").Value = "'6.435"
$ws.Range("E41").Value = "  -2.47%  "
$ws.Range("D42").Value = "'0.4529"
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").Value = "'0.8413"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "'101.94"
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("D47").Value = "'9.831"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("D48").Value = "'7.450"
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("D49").Value = "'36.73"
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("D50").Value = "'0.06097"
$ws.Range("E50").Value = "  +2.39%  "
$ws.Range("D51").Value = "'0.4187"
$ws.Range("E51").Value = "  +1.18%  "
